$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; all existing rows from 28 downward shift to 29+
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with its data
$ws.Cells.Item(28, 1).Value = 10
$ws.Cells.Item(28, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(28, 3).Value = "La Araucanía"
$ws.Cells.Item(28, 4).Value = 45107
$ws.Cells.Item(28, 5).Value = 9
$ws.Cells.Item(28, 6).Value = 300000001
$ws.Cells.Item(28, 7).Value = "Rabanito"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 65
$ws.Cells.Item(28, 11).Value = 7000
$ws.Cells.Item(28, 12).Value = 7000
$ws.Cells.Item(28, 13).Value = 7000
$ws.Cells.Item(28, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(28, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(28, 16).Value = 583
$ws.Cells.Item(28, 17).Value = 12
$ws.Cells.Item(28, 18).Value = "Hortaliza"
